$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Preserve the cells existing style while forcing the assigned
    # value to be stored as text (avoids numeric auto-coercion of
    # strings that look like plain decimal numbers, e.g. "606.77").
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = '65.814.05'
$ws.Range("E2").Value = '  +0.94%  '
$ws.Range("D3").Value = '2.700.91'
$ws.Range("E3").Value = '  +1.82%  '
$ws.Range("E4").Value = '  +0.03%  '
Set-TextValue $ws.Range("D5") '606.77'
$ws.Range("E5").Value = '  +1.80%  '
Set-TextValue $ws.Range("D6") '158.02'
$ws.Range("E6").Value = '  +1.02%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -0.79%  '
$ws.Range("E9").Value = '  +4.95%  '
Set-TextValue $ws.Range("D10") '6.06'
$ws.Range("E10").Value = '  +4.50%  '
Set-TextValue $ws.Range("D11") '0.402'
$ws.Range("E11").Value = '  +0.48%  '
$ws.Range("E12").Value = '  +1.19%  '
Set-TextValue $ws.Range("D13") '30.09'
$ws.Range("E13").Value = '  +3.58%  '
Set-TextValue $ws.Range("D14") '0.0000203'
$ws.Range("E14").Value = '  +9.95%  '
$ws.Range("D15").Value = '3.188.12'
$ws.Range("E15").Value = '  +1.88%  '
$ws.Range("D16").Value = '65.687.94'
$ws.Range("E16").Value = '  +0.94%  '
$ws.Range("D17").Value = '2.703.73'
$ws.Range("E17").Value = '  +0.47%  '
Set-TextValue $ws.Range("D18") '12.76'
$ws.Range("E18").Value = '  +1.00%  '
Set-TextValue $ws.Range("D19") '4.87'
$ws.Range("E19").Value = '  +1.31%  '
Set-TextValue $ws.Range("D20") '360.17'
$ws.Range("E20").Value = '  +1.44%  '
Set-TextValue $ws.Range("D21") '7.54'
$ws.Range("E21").Value = '  +3.50%  '
Set-TextValue $ws.Range("D22") '0.999'
$ws.Range("E22").Value = '  -0.19%  '
Set-TextValue $ws.Range("D23") '70.22'
$ws.Range("E23").Value = '  +2.86%  '
Set-TextValue $ws.Range("D24") '9.83'
$ws.Range("E24").Value = '  +3.26%  '
$ws.Range("E25").Value = '  +11.84%  '
$ws.Range("E26").Value = '  -4.03%  '
$ws.Range("E27").Value = '  +3.20%  '
$ws.Range("E28").Value = '  +3.70%  '
Set-TextValue $ws.Range("D29") '8.31'
$ws.Range("E29").Value = '  +1.69%  '
Set-TextValue $ws.Range("D30") '2.19'
$ws.Range("E30").Value = '  +4.12%  '
Set-TextValue $ws.Range("D31") '1.00'
$ws.Range("E31").Value = '  +0.19%  '
Set-TextValue $ws.Range("D32") '535.81'
$ws.Range("E32").Value = '  +2.57%  '
Set-TextValue $ws.Range("D33") '1.79'
$ws.Range("E33").Value = '  +0.19%  '
Set-TextValue $ws.Range("D34") '6.72'
$ws.Range("E34").Value = '  +5.51%  '
Set-TextValue $ws.Range("D35") '5.46'
$ws.Range("E35").Value = '  -3.00%  '
$ws.Range("E36").Value = '  +0.82%  '
Set-TextValue $ws.Range("D37") '20.76'
$ws.Range("E37").Value = '  +2.24%  '
Set-TextValue $ws.Range("D38") '162.52'
$ws.Range("E38").Value = '  -1.55%  '
$ws.Range("E39").Value = '  -0.89%  '
Set-TextValue $ws.Range("D40") '0.999'
$ws.Range("E40").Value = '  -0.08%  '
Set-TextValue $ws.Range("D41") '1.00'
$ws.Range("E41").Value = '  +0.04%  '
Set-TextValue $ws.Range("D42") '42.80'
$ws.Range("E42").Value = '  +1.55%  '
Set-TextValue $ws.Range("D43") '168.15'
$ws.Range("E43").Value = '  +1.54%  '
Set-TextValue $ws.Range("D44") '4.18'
$ws.Range("E44").Value = '  +1.87%  '
Set-TextValue $ws.Range("D45") '0.0619'
$ws.Range("E45").Value = '  -0.08%  '
Set-TextValue $ws.Range("D46") '23.59'
$ws.Range("E46").Value = '  +2.27%  '
Set-TextValue $ws.Range("D47") '2.27'
$ws.Range("E47").Value = '  +2.51%  '
Set-TextValue $ws.Range("D48") '0.0267'
$ws.Range("E48").Value = '  +4.52%  '
Set-TextValue $ws.Range("D49") '0.660'
$ws.Range("E49").Value = '  +1.54%  '
$ws.Range("E50").Value = '  +7.75%  '
Set-TextValue $ws.Range("D51") '0.0986'
$ws.Range("E51").Value = '  -0.22%  '
